$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column (Price) values: force text entry (NumberFormat '@') so numeric-looking
# strings like '210.32' are stored as text, matching the source inlineStr cells,
# then restore the default style so no stray NumberFormat/quotePrefix persists.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.933.38'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.46%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.592.79'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.53%  '

$ws.Range("E4").Value = '  +0.30%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.32'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.53%  '

$ws.Range("E6").Value = '  +0.30%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.482'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.11%  '

$ws.Range("E8").Value = '  -1.08%  '

$ws.Range("E9").Value = '  -1.42%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '17.92'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.74%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0809'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.69%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.814.52'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.60%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.589.21'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.26%  '

$ws.Range("E14").Value = '  -0.91%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.512'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.42%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.933.46'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.55%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '59.96'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.56%  '

$ws.Range("E18").Value = '  -0.53%  '

$ws.Range("E19").Value = '  +0.21%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '199.45'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +3.92%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.22'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.67%  '

$ws.Range("E22").Value = '  -2.22%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.99'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.76%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.80'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.91%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '141.93'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.46%  '

$ws.Range("E26").Value = '  +0.27%  '

$ws.Range("E27").Value = '  -8.33%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.06'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.68%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.43'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.68%  '

$ws.Range("E30").Value = '  +0.22%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.0474'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.15%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.10'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.75%  '

$ws.Range("E33").Value = '  -2.79%  '

$ws.Range("E34").Value = '  -2.17%  '

$ws.Range("E35").Value = '  +1.94%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.123.01'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.63%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0161'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +7.01%  '

$ws.Range("E38").Value = '  +0.32%  '

$ws.Range("E39").Value = '  -0.59%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.782'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.25%  '

$ws.Range("E41").Value = '  -3.56%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.777'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.95%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.726.41'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.59%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '92.43'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.97%  '

$ws.Range("E45").Value = '  -1.50%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.48'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.26%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '53.11'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.38%  '

$ws.Range("E48").Value = '  -1.39%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.408'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.44%  '

$ws.Range("E50").Value = '  +0.56%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0₇0915'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -18.12%  '
